$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.275.82"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.351.01"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.67"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.45"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").Value = "  +3.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.582"
$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.23"
$ws.Range("E11").Value = "  +5.74%  "

$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "686.03"
$ws.Range("E13").Value = "  +3.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.894.99"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.318.32"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.343.14"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  +2.27%  "

$ws.Range("E21").Value = "  +0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.43"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.95"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.91"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.70"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.92"
$ws.Range("E28").Value = "  -1.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.50"
$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("E30").Value = "  -6.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "561.13"
$ws.Range("E31").Value = "  -4.91%  "

$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.82"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.697.59"
$ws.Range("E36").Value = "  -0.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.28"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("E38").Value = "  +4.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.65"
$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.336"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.83"
$ws.Range("E50").Value = "  +2.86%  "

$ws.Range("E51").Value = "  -0.77%  "
